$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.890575333333333
$ws.Range("H2").Value = 5.671726
$ws.Range("I2").Value = 0.006346320422088561
$ws.Range("J2").Value = 0.00634632042208856
$ws.Range("M2").Value = 77.08952333333333
$ws.Range("N2").Value = 231.26857
$ws.Range("O2").Value = 0.2403816673726824
$ws.Range("P2").Value = 0.2403816673726824
$ws.Range("Q2").Value = 145.7435512724244
$ws.Range("R2").Value = 1311.69196145182
$ws.Range("S2").Value = 0.001525539084742954
$ws.Range("T2").Value = 0.001525539084742954

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.890575333333333
$ws.Range("H3").Value = 5.671726
$ws.Range("I3").Value = 0.006346320422088561
$ws.Range("J3").Value = 0.00634632042208856
$ws.Range("O3").Value = 0.3167483425780597
$ws.Range("P3").Value = 0.3167483425780597
$ws.Range("Q3").Value = 192.0447129414791
$ws.Range("R3").Value = 1728.402416473312
$ws.Range("S3").Value = 0.002010186475165844
$ws.Range("T3").Value = 0.002010186475165843

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.890575333333333
$ws.Range("H4").Value = 5.671726
$ws.Range("I4").Value = 0.006346320422088561
$ws.Range("J4").Value = 0.00634632042208856
$ws.Range("O4").Value = 0.4428699900492579
$ws.Range("P4").Value = 0.4428699900492579
$ws.Range("Q4").Value = 268.5123445861298
$ws.Range("R4").Value = 2416.611101275168
$ws.Range("S4").Value = 0.002810594862179763
$ws.Range("T4").Value = 0.002810594862179763

$ws.Range("I5").Value = 0.8887896079640043
$ws.Range("J5").Value = 0.8887896079640044
$ws.Range("M5").Value = 77.08952333333333
$ws.Range("N5").Value = 231.26857
$ws.Range("O5").Value = 0.2403816673726824
$ws.Range("P5").Value = 0.2403816673726824
$ws.Range("Q5").Value = 20411.09575051524
$ws.Range("R5").Value = 183699.8617546372
$ws.Range("S5").Value = 0.2136487279059001
$ws.Range("T5").Value = 0.2136487279059001

$ws.Range("I6").Value = 0.8887896079640043
$ws.Range("J6").Value = 0.8887896079640044
$ws.Range("O6").Value = 0.3167483425780597
$ws.Range("P6").Value = 0.3167483425780597
$ws.Range("S6").Value = 0.2815226352232018
$ws.Range("T6").Value = 0.2815226352232018

$ws.Range("I7").Value = 0.8887896079640043
$ws.Range("J7").Value = 0.8887896079640044
$ws.Range("O7").Value = 0.4428699900492579
$ws.Range("P7").Value = 0.4428699900492579
$ws.Range("S7").Value = 0.3936182448349024
$ws.Range("T7").Value = 0.3936182448349025

$ws.Range("I8").Value = 0.104864071613907
$ws.Range("J8").Value = 0.104864071613907
$ws.Range("M8").Value = 77.08952333333333
$ws.Range("N8").Value = 231.26857
$ws.Range("O8").Value = 0.2403816673726824
$ws.Range("P8").Value = 0.2403816673726824
$ws.Range("Q8").Value = 2408.20840761566
$ws.Range("R8").Value = 21673.87566854094
$ws.Range("S8").Value = 0.02520740038203933
$ws.Range("T8").Value = 0.02520740038203933

$ws.Range("I9").Value = 0.104864071613907
$ws.Range("J9").Value = 0.104864071613907
$ws.Range("O9").Value = 0.3167483425780597
$ws.Range("P9").Value = 0.3167483425780597
$ws.Range("S9").Value = 0.03321552087969199
$ws.Range("T9").Value = 0.03321552087969198

$ws.Range("I10").Value = 0.104864071613907
$ws.Range("J10").Value = 0.104864071613907
$ws.Range("O10").Value = 0.4428699900492579
$ws.Range("P10").Value = 0.4428699900492579
$ws.Range("S10").Value = 0.04644115035217564
$ws.Range("T10").Value = 0.04644115035217564

